## TC32_Verify_store_location.xlsx — "Update TC32 : Logixal QA box 2 env"
##
## The Contact-Us test block on the "TC32_Verify_store_location" sheet was
## rearranged/trimmed (two rows removed, the email-clear/enter-email steps
## moved after the message-entry step), and a stray duplicated row was
## removed from the "Testdata" object-repository sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: TC32_Verify_store_location
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Drop two rows from the bottom of the sheet (30 -> 28 rows); the
# remaining rows 18-28 get rewritten below with the new step order.
$ws1.Rows("29:30").Delete()

# Row 18: Select "Other" from the contact-us reason dropdown
$ws1.Range("B18").Value = "DROPDOWN_SELECT"
$ws1.Range("C18").Value = "ContactUsDropdown"
$ws1.Range("D18").Value = "CSS"
$ws1.Range("E18").Value = "ContactUsDropdown"

# Row 19: Click into the message box before entering text
$ws1.Range("B19").Value = "CLICK_PRE_ENTERTEXT"
$ws1.Range("C19").Value = "ContactUsMessage"
$ws1.Range("D19").Value = "CSS"
$ws1.Range("E19").Value = ""

# Row 20: Enter the message text
$ws1.Range("B20").Value = "ENTERTEXT"
$ws1.Range("C20").Value = "ContactUsMessage"
$ws1.Range("D20").Value = "CSS"
$ws1.Range("E20").Value = "ContactUsMessage"

# Row 21: Clear the email field
$ws1.Range("B21").Value = "CLEAR_TEXT"
$ws1.Range("C21").Value = "ContactUsEmail"
$ws1.Range("D21").Value = "CSS"
$ws1.Range("E21").Value = ""

# Row 22: Enter the email address
$ws1.Range("B22").Value = "ENTERTEXT"
$ws1.Range("C22").Value = "ContactUsEmail"
$ws1.Range("D22").Value = "CSS"
$ws1.Range("E22").Value = "ContactUsEmail"

# Row 23: Click "Send"
$ws1.Range("B23").Value = "CLICK"
$ws1.Range("C23").Value = "SendContactusDetails"
$ws1.Range("D23").Value = "CSS"
$ws1.Range("E23").Value = ""

# Row 24: Wait
$ws1.Range("B24").Value = "WAIT"
$ws1.Range("C24").Value = ""
$ws1.Range("D24").Value = ""
$ws1.Range("E24").Value = ""

# Row 25: Wait
$ws1.Range("B25").Value = "WAIT"
$ws1.Range("C25").Value = ""
$ws1.Range("D25").Value = ""
$ws1.Range("E25").Value = ""

# Row 26: Verify the confirmation message
$ws1.Range("B26").Value = "VERIFY_WEBELEMENT_PRESENT"
$ws1.Range("C26").Value = "ContactUsConfirmation"
$ws1.Range("D26").Value = "CSS"
$ws1.Range("E26").Value = "ContactUsConfirmation"

# Row 27: Open "My account"
$ws1.Range("B27").Value = "CLICK"
$ws1.Range("C27").Value = "MyaccountSection"
$ws1.Range("D27").Value = "CSS"
$ws1.Range("E27").Value = ""

# Row 28: Logout
$ws1.Range("B28").Value = "CLICK"
$ws1.Range("C28").Value = "Logout"
$ws1.Range("D28").Value = "CSS"
$ws1.Range("E28").Value = ""

# Restore the sheet's view: scrolled near the top of the Contact-Us block
# with the new email-clear row selected.
$ws1.Select()
$ws1.Rows(21).Select()

# ---------------------------------------------------------------------
# Sheet 2: Testdata (object repository)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 16 duplicated row 13 (ContactUsConfirmation / TRUE) - drop it.
$ws2.Rows(16).Delete()

$ws2.Select()
$ws2.Range("B12:B15").Select()
